$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-26: colo location table was reordered (AAE moved to top, ORN moved earlier)
$ws.Cells.Item(4, 1).Value = 'AAE'
$ws.Cells.Item(4, 2).Value = 'Annaba, Algeria'
$ws.Cells.Item(4, 3).Value = 36.85596
$ws.Cells.Item(4, 4).Value = 7.79207
$ws.Cells.Item(4, 5).Value = 'DZ'
$ws.Cells.Item(4, 6).Value = 'Africa'
$ws.Cells.Item(4, 7).Value = 'Annaba'

$ws.Cells.Item(5, 1).Value = 'TNR'
$ws.Cells.Item(5, 2).Value = 'Antananarivo, Madagascar'
$ws.Cells.Item(5, 3).Value = -18.91368
$ws.Cells.Item(5, 4).Value = 47.53613
$ws.Cells.Item(5, 5).Value = 'MG'
$ws.Cells.Item(5, 6).Value = 'Africa'
$ws.Cells.Item(5, 7).Value = 'Antananarivo'

$ws.Cells.Item(6, 1).Value = 'CPT'
$ws.Cells.Item(6, 2).Value = 'Cape Town, South Africa'
$ws.Cells.Item(6, 3).Value = -33.9648017883
$ws.Cells.Item(6, 4).Value = 18.6016998291
$ws.Cells.Item(6, 5).Value = 'ZA'
$ws.Cells.Item(6, 6).Value = 'Africa'
$ws.Cells.Item(6, 7).Value = 'Cape Town'

$ws.Cells.Item(7, 1).Value = 'CMN'
$ws.Cells.Item(7, 2).Value = 'Casablanca, Morocco'
$ws.Cells.Item(7, 3).Value = 33.3675003052
$ws.Cells.Item(7, 4).Value = -7.5899701118
$ws.Cells.Item(7, 5).Value = 'MA'
$ws.Cells.Item(7, 6).Value = 'Africa'
$ws.Cells.Item(7, 7).Value = 'Casablanca'

$ws.Cells.Item(8, 1).Value = 'DKR'
$ws.Cells.Item(8, 2).Value = 'Dakar, Senegal'
$ws.Cells.Item(8, 3).Value = 14.7412099
$ws.Cells.Item(8, 4).Value = -17.4889771
$ws.Cells.Item(8, 5).Value = 'SN'
$ws.Cells.Item(8, 6).Value = 'Africa'
$ws.Cells.Item(8, 7).Value = 'Dakar'

$ws.Cells.Item(9, 1).Value = 'DAR'
$ws.Cells.Item(9, 2).Value = 'Dar Es Salaam, Tanzania'
$ws.Cells.Item(9, 3).Value = -6.8781099319
$ws.Cells.Item(9, 4).Value = 39.2025985718
$ws.Cells.Item(9, 5).Value = 'TZ'
$ws.Cells.Item(9, 6).Value = 'Africa'
$ws.Cells.Item(9, 7).Value = 'Dar es Salaam'

$ws.Cells.Item(10, 1).Value = 'JIB'
$ws.Cells.Item(10, 2).Value = 'Djibouti City, Djibouti'
$ws.Cells.Item(10, 3).Value = 11.5473003387
$ws.Cells.Item(10, 4).Value = 43.1595001221
$ws.Cells.Item(10, 5).Value = 'DJ'
$ws.Cells.Item(10, 6).Value = 'Africa'
$ws.Cells.Item(10, 7).Value = 'Djibouti'

$ws.Cells.Item(11, 1).Value = 'DUR'
$ws.Cells.Item(11, 2).Value = 'Durban, South Africa'
$ws.Cells.Item(11, 3).Value = -29.6144444444
$ws.Cells.Item(11, 4).Value = 31.1197222222
$ws.Cells.Item(11, 5).Value = 'ZA'
$ws.Cells.Item(11, 6).Value = 'Africa'
$ws.Cells.Item(11, 7).Value = 'Durban'

$ws.Cells.Item(12, 1).Value = 'GBE'
$ws.Cells.Item(12, 2).Value = 'Gaborone, Botswana'
$ws.Cells.Item(12, 3).Value = -24.6282
$ws.Cells.Item(12, 4).Value = 25.9231
$ws.Cells.Item(12, 5).Value = 'BW'
$ws.Cells.Item(12, 6).Value = 'Africa'
$ws.Cells.Item(12, 7).Value = 'Gaborone'

$ws.Cells.Item(13, 1).Value = 'HRE'
$ws.Cells.Item(13, 2).Value = 'Harare, Zimbabwe'
$ws.Cells.Item(13, 3).Value = -17.9318008423
$ws.Cells.Item(13, 4).Value = 31.0928001404
$ws.Cells.Item(13, 5).Value = 'ZW'
$ws.Cells.Item(13, 6).Value = 'Africa'
$ws.Cells.Item(13, 7).Value = 'Harare'

$ws.Cells.Item(14, 1).Value = 'JNB'
$ws.Cells.Item(14, 2).Value = 'Johannesburg, South Africa'
$ws.Cells.Item(14, 3).Value = -26.133333
$ws.Cells.Item(14, 4).Value = 28.25
$ws.Cells.Item(14, 5).Value = 'ZA'
$ws.Cells.Item(14, 6).Value = 'Africa'
$ws.Cells.Item(14, 7).Value = 'Johannesburg'

$ws.Cells.Item(15, 1).Value = 'KGL'
$ws.Cells.Item(15, 2).Value = 'Kigali, Rwanda'
$ws.Cells.Item(15, 3).Value = -1.9686299563
$ws.Cells.Item(15, 4).Value = 30.1394996643
$ws.Cells.Item(15, 5).Value = 'RW'
$ws.Cells.Item(15, 6).Value = 'Africa'
$ws.Cells.Item(15, 7).Value = 'Kigali'

$ws.Cells.Item(16, 1).Value = 'LOS'
$ws.Cells.Item(16, 2).Value = 'Lagos, Nigeria'
$ws.Cells.Item(16, 3).Value = 6.5773701668
$ws.Cells.Item(16, 4).Value = 3.321160078
$ws.Cells.Item(16, 5).Value = 'NG'
$ws.Cells.Item(16, 6).Value = 'Africa'
$ws.Cells.Item(16, 7).Value = 'Lagos'

$ws.Cells.Item(17, 1).Value = 'LAD'
$ws.Cells.Item(17, 2).Value = 'Luanda, Angola'
$ws.Cells.Item(17, 3).Value = -8.858369827300001
$ws.Cells.Item(17, 4).Value = 13.2312002182
$ws.Cells.Item(17, 5).Value = 'AO'
$ws.Cells.Item(17, 6).Value = 'Africa'
$ws.Cells.Item(17, 7).Value = 'Luanda'

$ws.Cells.Item(18, 1).Value = 'MPM'
$ws.Cells.Item(18, 2).Value = 'Maputo, Mozambique'
$ws.Cells.Item(18, 3).Value = -25.9207992554
$ws.Cells.Item(18, 4).Value = 32.5726013184
$ws.Cells.Item(18, 5).Value = 'MZ'
$ws.Cells.Item(18, 6).Value = 'Africa'
$ws.Cells.Item(18, 7).Value = 'Maputo'

$ws.Cells.Item(19, 1).Value = 'MBA'
$ws.Cells.Item(19, 2).Value = 'Mombasa, Kenya'
$ws.Cells.Item(19, 3).Value = -4.0348300934
$ws.Cells.Item(19, 4).Value = 39.5942001343
$ws.Cells.Item(19, 5).Value = 'KE'
$ws.Cells.Item(19, 6).Value = 'Africa'
$ws.Cells.Item(19, 7).Value = 'Mombasa'

$ws.Cells.Item(20, 1).Value = 'NBO'
$ws.Cells.Item(20, 2).Value = 'Nairobi, Kenya'
$ws.Cells.Item(20, 3).Value = -1.319239974
$ws.Cells.Item(20, 4).Value = 36.9277992249
$ws.Cells.Item(20, 5).Value = 'KE'
$ws.Cells.Item(20, 6).Value = 'Africa'
$ws.Cells.Item(20, 7).Value = 'Nairobi'

$ws.Cells.Item(21, 1).Value = 'ORN'
$ws.Cells.Item(21, 2).Value = 'Oran, Algeria'
$ws.Cells.Item(21, 3).Value = 35.6911
$ws.Cells.Item(21, 4).Value = -0.6415999999999999
$ws.Cells.Item(21, 5).Value = 'DZ'
$ws.Cells.Item(21, 6).Value = 'Africa'
$ws.Cells.Item(21, 7).Value = 'Oran'

$ws.Cells.Item(22, 1).Value = 'OUA'
$ws.Cells.Item(22, 2).Value = 'Ouagadougou, Burkina Faso'
$ws.Cells.Item(22, 3).Value = 12.3531999588
$ws.Cells.Item(22, 4).Value = -1.5124200583
$ws.Cells.Item(22, 5).Value = 'BF'
$ws.Cells.Item(22, 6).Value = 'Africa'
$ws.Cells.Item(22, 7).Value = 'Ouagadougou'

$ws.Cells.Item(23, 1).Value = 'MRU'
$ws.Cells.Item(23, 2).Value = 'Port Louis, Mauritius'
$ws.Cells.Item(23, 3).Value = -20.4302005768
$ws.Cells.Item(23, 4).Value = 57.6836013794
$ws.Cells.Item(23, 5).Value = 'MU'
$ws.Cells.Item(23, 6).Value = 'Africa'
$ws.Cells.Item(23, 7).Value = 'Port Louis'

$ws.Cells.Item(24, 1).Value = 'RUN'
$ws.Cells.Item(24, 2).Value = 'Réunion, France'
$ws.Cells.Item(24, 3).Value = -20.8871002197
$ws.Cells.Item(24, 4).Value = 55.5102996826
$ws.Cells.Item(24, 5).Value = 'RE'
$ws.Cells.Item(24, 6).Value = 'Africa'
$ws.Cells.Item(24, 7).Value = 'Saint-Denis'

$ws.Cells.Item(25, 1).Value = 'TUN'
$ws.Cells.Item(25, 2).Value = 'Tunis, Tunisia'
$ws.Cells.Item(25, 3).Value = 36.8510017395
$ws.Cells.Item(25, 4).Value = 10.2271995544
$ws.Cells.Item(25, 5).Value = 'TN'
$ws.Cells.Item(25, 6).Value = 'Africa'
$ws.Cells.Item(25, 7).Value = 'Tunis'

$ws.Cells.Item(26, 1).Value = 'FIH'
$ws.Cells.Item(26, 2).Value = 'Kinshasa, DR Congo'
$ws.Cells.Item(26, 3).Value = -4.3857498169
$ws.Cells.Item(26, 4).Value = 15.4446001053
$ws.Cells.Item(26, 5).Value = 'CD'
$ws.Cells.Item(26, 6).Value = 'Africa'
$ws.Cells.Item(26, 7).Value = 'Kinshasa'

# Row 109: city/country name updated
$ws.Cells.Item(109, 2).Value = 'Shenzhen, China'
